$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.041545723200529
$ws.Range("D2").Value = 1.039297054728629
$ws.Range("E2").Value = 1.048715094248676
$ws.Range("F2").Value = 1.057042820102352
$ws.Range("I2").Value = 1.038489970946423
$ws.Range("J2").Value = 1.046626351535078
$ws.Range("K2").Value = 1.042082643582883
$ws.Range("L2").Value = 1.051474166906976
$ws.Range("M2").Value = 1.059778872817172
$ws.Range("N2").Value = 1.019363311645429

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.042714220149611
$ws.Range("D3").Value = 1.039879405151342
$ws.Range("E3").Value = 1.049791665577416
$ws.Range("F3").Value = 1.058290854848098
$ws.Range("I3").Value = 1.038721569036194
$ws.Range("J3").Value = 1.04743983060988
$ws.Range("K3").Value = 1.042475784736046
$ws.Range("L3").Value = 1.0523621553137
$ws.Range("M3").Value = 1.060839562023407
$ws.Range("N3").Value = 1.019639937002866

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.043470160793303
$ws.Range("D4").Value = 1.040256136789622
$ws.Range("E4").Value = 1.050488512241968
$ws.Range("F4").Value = 1.059098953709223
$ws.Range("I4").Value = 1.038870149398911
$ws.Range("J4").Value = 1.047965522699627
$ws.Range("K4").Value = 1.042729385589925
$ws.Range("L4").Value = 1.052936380819574
$ws.Range("M4").Value = 1.061525873563215
$ws.Range("N4").Value = 1.019818538618149

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.043787922607966
$ws.Range("D5").Value = 1.040414492663921
$ws.Range("E5").Value = 1.050781523167494
$ws.Range("F5").Value = 1.059438808208897
$ws.Range("I5").Value = 1.038932306287537
$ws.Range("J5").Value = 1.048186360797449
$ws.Range("K5").Value = 1.042835810251191
$ws.Range("L5").Value = 1.053177699229835
$ws.Range("M5").Value = 1.061814393873732
$ws.Range("N5").Value = 1.019893528663475

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.043841274168446
$ws.Range("D6").Value = 1.040441079985529
$ws.Range("E6").Value = 1.050830724332821
$ws.Range("F6").Value = 1.05949587891759
$ws.Range("I6").Value = 1.038942724738982
$ws.Range("J6").Value = 1.048223430956712
$ws.Range("K6").Value = 1.042853668313035
$ws.Range("L6").Value = 1.053218212632235
$ws.Range("M6").Value = 1.061862837435651
$ws.Range("N6").Value = 1.019906114309239

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.043474406880377
$ws.Range("D7").Value = 1.04025825283674
$ws.Range("E7").Value = 1.050492427244944
$ws.Range("F7").Value = 1.059103494350397
$ws.Range("I7").Value = 1.038870981145933
$ws.Range("J7").Value = 1.04796847418812
$ws.Range("K7").Value = 1.042730808385334
$ws.Range("L7").Value = 1.052939605664927
$ws.Range("M7").Value = 1.061529728804364
$ws.Range("N7").Value = 1.019819541008292

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.041940655203081
$ws.Range("D8").Value = 1.03949388038797
$ws.Range("E8").Value = 1.049078878802345
$ws.Range("F8").Value = 1.057464488834086
$ws.Range("I8").Value = 1.038568505401312
$ws.Range("J8").Value = 1.046901412505239
$ws.Range("K8").Value = 1.042215670260238
$ws.Range("L8").Value = 1.051774342235743
$ws.Range("M8").Value = 1.060137343539694
$ws.Range("N8").Value = 1.019456880032092

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.039236735051743
$ws.Range("D9").Value = 1.038146329494466
$ws.Range("E9").Value = 1.04658977426787
$ws.Range("F9").Value = 1.054580410009486
$ws.Range("I9").Value = 1.038025711455913
$ws.Range("J9").Value = 1.045015847904776
$ws.Range("K9").Value = 1.041301921529563
$ws.Range("L9").Value = 1.049718189457522
$ws.Range("M9").Value = 1.057683544045121
$ws.Range("N9").Value = 1.018814805442793

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.0374331833529
$ws.Range("D10").Value = 1.03724759907259
$ws.Range("E10").Value = 1.044931493394102
$ws.Range("F10").Value = 1.052660341371929
$ws.Range("I10").Value = 1.037657267663314
$ws.Range("J10").Value = 1.043755216737393
$ws.Range("K10").Value = 1.040688746036693
$ws.Range("L10").Value = 1.048345481394898
$ws.Range("M10").Value = 1.056047451840392
$ws.Range("N10").Value = 1.018384716141186

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.036651981049461
$ws.Range("D11").Value = 1.036858364958461
$ws.Range("E11").Value = 1.044213692908706
$ws.Range("F11").Value = 1.051829540958568
$ws.Range("I11").Value = 1.037496166704893
$ws.Range("J11").Value = 1.043208487107807
$ws.Range("K11").Value = 1.040422289326136
$ws.Range("L11").Value = 1.047750611652729
$ws.Range("M11").Value = 1.055338936956213
$ws.Range("N11").Value = 1.018197996559542

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.036361767802503
$ws.Range("D12").Value = 1.036713775316047
$ws.Range("E12").Value = 1.043947105624576
$ws.Range("F12").Value = 1.051521033470534
$ws.Range("I12").Value = 1.037436091863632
$ws.Range("J12").Value = 1.043005275896995
$ws.Range("K12").Value = 1.040323173498792
$ws.Range("L12").Value = 1.047529577592216
$ws.Range("M12").Value = 1.055075750122164
$ws.Range("N12").Value = 1.018128567017202

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.036424021367965
$ws.Range("D13").Value = 1.036744790757094
$ws.Range("E13").Value = 1.044004287902231
$ws.Range("F13").Value = 1.051587205383369
$ws.Range("I13").Value = 1.03744898874568
$ws.Range("J13").Value = 1.043048871358694
$ws.Range("K13").Value = 1.040344440596385
$ws.Range("L13").Value = 1.047576993467803
$ws.Range("M13").Value = 1.05513220519617
$ws.Range("N13").Value = 1.018143463223191

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.036627992757985
$ws.Range("D14").Value = 1.036846413347625
$ws.Range("E14").Value = 1.04419165599773
$ws.Range("F14").Value = 1.051804037831975
$ws.Range("I14").Value = 1.037491205689319
$ws.Range("J14").Value = 1.043191692294167
$ws.Range("K14").Value = 1.040414099273975
$ws.Range("L14").Value = 1.047732342403035
$ws.Range("M14").Value = 1.055317182120939
$ws.Range("N14").Value = 1.018192258994282

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.036753660912795
$ws.Range("D15").Value = 1.036909024981663
$ws.Range("E15").Value = 1.044307104370227
$ws.Range("F15").Value = 1.051937647173939
$ws.Range("I15").Value = 1.037517185830673
$ws.Range("J15").Value = 1.043279671531231
$ws.Range("K15").Value = 1.040456999488677
$ws.Range("L15").Value = 1.047828048308733
$ws.Range("M15").Value = 1.055431150762881
$ws.Range("N15").Value = 1.018222313922328

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.037485023686854
$ws.Range("D16").Value = 1.037273429698052
$ws.Range("E16").Value = 1.044979136530998
$ws.Range("F16").Value = 1.052715491400231
$ws.Range("I16").Value = 1.037667926482333
$ws.Range("J16").Value = 1.043791483002619
$ws.Range("K16").Value = 1.040706409961941
$ws.Range("L16").Value = 1.048384950771811
$ws.Range("M16").Value = 1.05609447191427
$ws.Range("N16").Value = 1.018397097800533

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.037943719222868
$ws.Range("D17").Value = 1.037501990958848
$ws.Range("E17").Value = 1.045400749590692
$ws.Range("F17").Value = 1.053203572526448
$ws.Range("I17").Value = 1.0377620637901
$ws.Range("J17").Value = 1.04411229566759
$ws.Range("K17").Value = 1.040862605152428
$ws.Range("L17").Value = 1.048734152182712
$ws.Range("M17").Value = 1.056510534357859
$ws.Range("N17").Value = 1.018506604263307

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.038211244338029
$ws.Range("D18").Value = 1.037635299295232
$ws.Range("E18").Value = 1.045646693178317
$ws.Range("F18").Value = 1.053488320183494
$ws.Range("I18").Value = 1.037816821746543
$ws.Range("J18").Value = 1.044299336588635
$ws.Range("K18").Value = 1.040953619698429
$ws.Range("L18").Value = 1.04893778940419
$ws.Range("M18").Value = 1.056753209288537
$ws.Range("N18").Value = 1.018570430454839

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.038302459439437
$ws.Range("D19").Value = 1.037680752678407
$ws.Range("E19").Value = 1.04573055766301
$ws.Range("F19").Value = 1.053585421697241
$ws.Range("I19").Value = 1.037835467227415
$ws.Range("J19").Value = 1.044363098548403
$ws.Range("K19").Value = 1.040984637779021
$ws.Range("L19").Value = 1.049007216656929
$ws.Range("M19").Value = 1.056835953981364
$ws.Range("N19").Value = 1.018592185559433

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.037894508053287
$ws.Range("D20").Value = 1.037477469275817
$ws.Range("E20").Value = 1.045355512043778
$ws.Range("F20").Value = 1.053151200033265
$ws.Range("I20").Value = 1.037751979340152
$ws.Range("J20").Value = 1.04407788416763
$ws.Range("K20").Value = 1.04084585634425
$ws.Range("L20").Value = 1.048696690947583
$ws.Range("M20").Value = 1.056465895555505
$ws.Range("N20").Value = 1.01849486012718

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.036567929385696
$ws.Range("D21").Value = 1.03681648831713
$ws.Range("E21").Value = 1.044136479785637
$ws.Range("F21").Value = 1.051740183642961
$ws.Range("I21").Value = 1.037478780332874
$ws.Range("J21").Value = 1.043149638726605
$ws.Range("K21").Value = 1.040393590444406
$ws.Range("L21").Value = 1.04768659803644
$ws.Range("M21").Value = 1.055262711393397
$ws.Range("N21").Value = 1.018177891888958

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.035733624724877
$ws.Range("D22").Value = 1.036400841551134
$ws.Range("E22").Value = 1.04337023221851
$ws.Range("F22").Value = 1.050853534445661
$ws.Range("I22").Value = 1.037305650930912
$ws.Range("J22").Value = 1.042565252720941
$ws.Range("K22").Value = 1.040108411908454
$ws.Range("L22").Value = 1.047051089898612
$ws.Range("M22").Value = 1.054506146652585
$ws.Range("N22").Value = 1.017978175564896

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.036175928067674
$ws.Range("D23").Value = 1.036621189311366
$ws.Range("E23").Value = 1.043776415387335
$ws.Range("F23").Value = 1.051323515908884
$ws.Range("I23").Value = 1.037397558873025
$ws.Range("J23").Value = 1.042875119317025
$ws.Range("K23").Value = 1.040259668124242
$ws.Range("L23").Value = 1.047388025223544
$ws.Range("M23").Value = 1.054907223529232
$ws.Range("N23").Value = 1.018084089412463

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.037916744536877
$ws.Range("D24").Value = 1.037488549592338
$ws.Range("E24").Value = 1.045375952870015
$ws.Range("F24").Value = 1.053174864728299
$ws.Range("I24").Value = 1.037756536534698
$ws.Range("J24").Value = 1.044093433501447
$ws.Range("K24").Value = 1.040853424692327
$ws.Range("L24").Value = 1.048713618209181
$ws.Range("M24").Value = 1.056486065930823
$ws.Range("N24").Value = 1.018500166942133

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.039935921607596
$ws.Range("D25").Value = 1.038494772055206
$ws.Range("E25").Value = 1.047233065743814
$ws.Range("F25").Value = 1.055325540611598
$ws.Range("I25").Value = 1.038167196900261
$ws.Range("J25").Value = 1.045503940464386
$ws.Range("K25").Value = 1.041538856449429
$ws.Range("L25").Value = 1.050250092389391
$ws.Range("M25").Value = 1.058317945047652
$ws.Range("N25").Value = 1.019363311645429

Write-Output "vm_pu values updated for case with 380 kV"